$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-driving data: rows 2-10, columns A-T
# representing the full 3x3 Sending x Target cluster grid (ECs/FAPs/sCs)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf1"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.161357666666667
$ws.Range("H2").Value = 3.484073
$ws.Range("I2").Value = 0.1270850363824361
$ws.Range("J2").Value = 0.1270850363824361
$ws.Range("K2").Value = 2.0
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.619953333333334
$ws.Range("N2").Value = 7.85986
$ws.Range("O2").Value = 0.6100029538328192
$ws.Range("P2").Value = 0.6100029538328192
$ws.Range("Q2").Value = 3.042702889975556
$ws.Range("R2").Value = 27.38432600978
$ws.Range("S2").Value = 0.07752224758123731
$ws.Range("T2").Value = 0.07752224758123731

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf1"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.161357666666667
$ws.Range("H3").Value = 3.484073
$ws.Range("I3").Value = 0.1270850363824361
$ws.Range("J3").Value = 0.1270850363824361
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 0.5698483333333333
$ws.Range("N3").Value = 1.709545
$ws.Range("O3").Value = 0.1326776176306101
$ws.Range("P3").Value = 0.1326776176306101
$ws.Range("Q3").Value = 0.6617977307538889
$ws.Range("R3").Value = 5.956179576785
$ws.Range("S3").Value = 0.01686133986372102
$ws.Range("T3").Value = 0.01686133986372102

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf1"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.161357666666667
$ws.Range("H4").Value = 3.484073
$ws.Range("I4").Value = 0.1270850363824361
$ws.Range("J4").Value = 0.1270850363824361
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 1.105183
$ws.Range("N4").Value = 3.315549
$ws.Range("O4").Value = 0.2573194285365706
$ws.Range("P4").Value = 0.2573194285365706
$ws.Range("Q4").Value = 1.283512750119667
$ws.Range("R4").Value = 11.551614751077
$ws.Range("S4").Value = 0.03270144893747774
$ws.Range("T4").Value = 0.03270144893747774

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf1"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 5.168173666666667
$ws.Range("H5").Value = 15.504521
$ws.Range("I5").Value = 0.565542861868062
$ws.Range("J5").Value = 0.565542861868062
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.619953333333334
$ws.Range("N5").Value = 7.85986
$ws.Range("O5").Value = 0.6100029538328192
$ws.Range("P5").Value = 0.6100029538328192
$ws.Range("Q5").Value = 13.54037382522889
$ws.Range("R5").Value = 121.86336442706
$ws.Range("S5").Value = 0.3449828162585839
$ws.Range("T5").Value = 0.3449828162585839

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf1"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 5.168173666666667
$ws.Range("H6").Value = 15.504521
$ws.Range("I6").Value = 0.565542861868062
$ws.Range("J6").Value = 0.565542861868062
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.5698483333333333
$ws.Range("N6").Value = 1.709545
$ws.Range("O6").Value = 0.1326776176306101
$ws.Range("P6").Value = 0.1326776176306101
$ws.Range("Q6").Value = 2.945075150327222
$ws.Range("R6").Value = 26.505676352945
$ws.Range("S6").Value = 0.07503487958065165
$ws.Range("T6").Value = 0.07503487958065165

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf1"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 5.168173666666667
$ws.Range("H7").Value = 15.504521
$ws.Range("I7").Value = 0.565542861868062
$ws.Range("J7").Value = 0.565542861868062
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 1.105183
$ws.Range("N7").Value = 3.315549
$ws.Range("O7").Value = 0.2573194285365706
$ws.Range("P7").Value = 0.2573194285365706
$ws.Range("Q7").Value = 5.711777677447667
$ws.Range("R7").Value = 51.405999097029
$ws.Range("S7").Value = 0.1455251660288264
$ws.Range("T7").Value = 0.1455251660288264

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf1"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 2.808898333333333
$ws.Range("H8").Value = 8.426695
$ws.Range("I8").Value = 0.3073721017495019
$ws.Range("J8").Value = 0.3073721017495019
$ws.Range("K8").Value = 2.0
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.619953333333334
$ws.Range("N8").Value = 7.85986
$ws.Range("O8").Value = 0.6100029538328192
$ws.Range("P8").Value = 0.6100029538328192
$ws.Range("Q8").Value = 7.359182551411112
$ws.Range("R8").Value = 66.23264296270001
$ws.Range("S8").Value = 0.187497889992998
$ws.Range("T8").Value = 0.187497889992998

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf1"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 2.808898333333333
$ws.Range("H9").Value = 8.426695
$ws.Range("I9").Value = 0.3073721017495019
$ws.Range("J9").Value = 0.3073721017495019
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.5698483333333333
$ws.Range("N9").Value = 1.709545
$ws.Range("O9").Value = 0.1326776176306101
$ws.Range("P9").Value = 0.1326776176306101
$ws.Range("Q9").Value = 1.600646033752778
$ws.Range("R9").Value = 14.405814303775
$ws.Range("S9").Value = 0.04078139818623738
$ws.Range("T9").Value = 0.04078139818623738

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf1"
$ws.Range("C10").Value = "Fgfr3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 2.808898333333333
$ws.Range("H10").Value = 8.426695
$ws.Range("I10").Value = 0.3073721017495019
$ws.Range("J10").Value = 0.3073721017495019
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 1.105183
$ws.Range("N10").Value = 3.315549
$ws.Range("O10").Value = 0.2573194285365706
$ws.Range("P10").Value = 0.2573194285365706
$ws.Range("Q10").Value = 3.104346686728334
$ws.Range("R10").Value = 27.939120180555
$ws.Range("S10").Value = 0.07909281357026647
$ws.Range("T10").Value = 0.07909281357026647

Write-Output "Edit complete"
